# Append: 2025-12-21 06:27 JST
# New scrape pass replaces the previous 5 newest rows (rows 2-6) with fresh
# listings and drops the two oldest rows (7-8) that fell off the window.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("ランサーズ")

# --- Update rows 2-6 with the new scrape results -------------------------

$ws.Range("A2").Value = "2025-12-21 06:27:45"
$ws.Range("B2").Value = "【急募】ECサイトの自動購入Bot作成をお願いします。"
$ws.Range("D2").Value = "100,000 円 ~ 200,000 円 / 固定"
$ws.Range("F2").Value = "https://www.lancers.jp/work/detail/5458190"
$ws.Range("G2").Value = 143
$ws.Range("H2").Value = "★bot ◇サイト"

$ws.Range("A3").Value = "2025-12-21 06:27:45"
$ws.Range("B3").Value = "【受注メールを元にECサイト自動仕入ツール】"
$ws.Range("D3").Value = "50,000 円 ~ 100,000 円 / 固定"
$ws.Range("F3").Value = "https://www.lancers.jp/work/detail/5458166"
$ws.Range("G3").Value = 98
$ws.Range("H3").Value = "◆ツール ◇サイト"

$ws.Range("A4").Value = "2025-12-21 06:27:45"
$ws.Range("B4").Value = "【自動運転プロジェクト経験者募集】実証実験・開発を推進するプロジェクトマネージャー"
$ws.Range("D4").Value = "200,000 円 ~ 300,000 円 / 固定"
$ws.Range("F4").Value = "https://www.lancers.jp/work/detail/5431107"
$ws.Range("G4").Value = 68
$ws.Range("H4").Value = "◆開発"

$ws.Range("A5").Value = "2025-12-21 06:27:45"
$ws.Range("B5").Value = "初回 【急募】ECサイトの要件定義や基本設計ができる方を募集(1人月、フルリモート可、2025年12月〜)"
$ws.Range("D5").Value = "300,000 円 ~ 500,000 円 / 固定"
$ws.Range("F5").Value = "https://www.lancers.jp/work/detail/5425629"
$ws.Range("G5").Value = 45
$ws.Range("H5").Value = "◇サイト"

$ws.Range("A6").Value = "2025-12-21 06:27:45"
$ws.Range("B6").Value = "【急募】Notionでの社内向けダッシュボード作成依頼"
$ws.Range("D6").Value = "50,000 円 ~ 100,000 円 / 固定"
$ws.Range("F6").Value = "https://www.lancers.jp/work/detail/5458234"
$ws.Range("G6").Value = 18
$ws.Range("H6").ClearContents()

# --- Drop the two oldest rows (7 & 8) -------------------------------------

$ws.Range("A7:H8").EntireRow.Delete()

# --- Rebuild hyperlinks so F2:F6 point at the refreshed URLs --------------
# (stale F7/F8 hyperlink entries must go too, and the per-item Delete() is a
# no-op in this host, so wipe the sheet's hyperlinks and re-add the live ones)

$ws.Hyperlinks.Delete()
$ws.Hyperlinks.Add($ws.Range("F2"), "https://www.lancers.jp/work/detail/5458190")
$ws.Hyperlinks.Add($ws.Range("F3"), "https://www.lancers.jp/work/detail/5458166")
$ws.Hyperlinks.Add($ws.Range("F4"), "https://www.lancers.jp/work/detail/5431107")
$ws.Hyperlinks.Add($ws.Range("F5"), "https://www.lancers.jp/work/detail/5425629")
$ws.Hyperlinks.Add($ws.Range("F6"), "https://www.lancers.jp/work/detail/5458234")

# --- Column width tweaks (D: 30 -> 28, H: 25 -> 12) -----------------------
# ColumnWidth is in Excel character units, which round-trip to the saved XML
# width with a +0.8333... padding offset, so back that off here.

$ws.Columns.Item(4).ColumnWidth = 28 - 0.8333333333333
$ws.Columns.Item(8).ColumnWidth = 12 - 0.8333333333333
